# Adds three new rows (109, 110, 111) to the "BalancoResumido" sheet, each
# duplicating the contents of row 108 (columns A:Q). Row 108's stray empty
# column-R cell is removed, and the same stray empty column-R cell is
# recreated on the new last row (111).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BalancoResumido")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")

function Copy-RowValues($srcRow, $dstRow) {
    foreach ($col in $cols) {
        $srcCell = $ws.Range($col + $srcRow)
        $dstCell = $ws.Range($col + $dstRow)
        $val = $srcCell.Value()
        $dstCell.NumberFormat = "@"
        $dstCell.Value = $val
        $dstCell.Style = "Normal"
    }
}

Copy-RowValues 108 109
Copy-RowValues 108 110
Copy-RowValues 108 111

# Row 108 originally carried a trailing (empty) column-R cell; move it down
# to the new final row (111) instead.
$ws.Range("R108").ClearContents()
$ws.Range("R111").NumberFormat = "@"
$ws.Range("R111").Value = "X"
$ws.Range("R111").Value = ""
$ws.Range("R111").Style = "Normal"
